# ForecastResult2017.xlsx — "Can Forecast Room Nights Sold"
#
# - Adds two workbook-level defined names ("rns" / "Timeline") used by the
#   new FORECAST.ETS formula.
# - Re-labels the headers (drops the old "Forecast"/"Lower/Upper Confidence
#   Bound" columns, adds "Forecasted Date" + "Forecast").
# - Replaces the 24 rows of hard-coded text month-labels in column B with
#   live DATE() formulas for the most recent 6 months, and refreshes the
#   "Room Nights Sold" figures in column C to match.
# - Adds the forecast helper cells E2 (next month, via EOMONTH) and F2
#   (FORECAST.ETS over the new B/C history, using the new names).
# - Removes the old scratch FORECAST.ETS formulas that lived in D26:F26.
# - Grows the sheet down to row 362 and widens column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Defined names consumed by the new forecast formula in F2.
# ---------------------------------------------------------------------
$wb.Names.Add("rns", "='Sheet1'!`$C`$2:`$C`$7")
$wb.Names.Add("Timeline", "='Sheet1'!`$B`$2:`$B`$7")

# ---------------------------------------------------------------------
# 2. Header row (row 1): drop the old D1 "Forecast" column entirely,
#    relabel E1/F1.
# ---------------------------------------------------------------------
$ws.Range("D1").Clear()
$ws.Range("E1").Clear()
$ws.Range("E1").Value = "Forecasted Date"
$ws.Range("F1").Value = "Forecast"

# ---------------------------------------------------------------------
# 3. Column B/C history: wipe the old 24 text-date rows + values, then
#    lay down the 6 rows of live data the new forecast is based on.
# ---------------------------------------------------------------------
$ws.Range("B2:C25").ClearContents()

$ws.Range("B2").Formula = "=DATE(2016,07,31)"
$ws.Range("C2").Value = 1100

$ws.Range("B3").Formula = "=DATE(2016,08,31)"
$ws.Range("C3").Value = 1003

$ws.Range("B4").Formula = "=DATE(2016,09,30)"
$ws.Range("C4").Value = 1223

$ws.Range("B5").Formula = "=DATE(2016,10,31)"
$ws.Range("C5").Value = 1440

$ws.Range("B6").Formula = "=DATE(2016,11,30)"
$ws.Range("C6").Value = 1161

$ws.Range("B7").Formula = "=DATE(2016,12,31)"
$ws.Range("C7").Value = 2686

# ---------------------------------------------------------------------
# 4. New forecast helper cells on row 2.
# ---------------------------------------------------------------------
$ws.Range("E2").Formula = "=EOMONTH(DATE(2016,12,31),1)"
$ws.Range("E2").Style = $ws.Range("B2").Style

$ws.Range("F2").Formula = '=_xlfn.FORECAST.ETS($E2,rns,Timeline,1,1)'
$ws.Range("F2").Style = $ws.Range("F1").Style

$ws.Range("G2").Style = $ws.Range("F1").Style

# ---------------------------------------------------------------------
# 5. Remove the old scratch FORECAST.ETS formulas in row 26.
# ---------------------------------------------------------------------
$ws.Range("D26:F26").ClearContents()

# ---------------------------------------------------------------------
# 6. Grow the sheet down to row 362 (27 new blank, styled rows) and
#    widen column G.
# ---------------------------------------------------------------------
$ws.Range("B336:B362").Style = $ws.Range("B335").Style

$ws.Columns.Item(7).ColumnWidth = 22.6

# ---------------------------------------------------------------------
# 7. Selection moves to F2.
# ---------------------------------------------------------------------
$ws.Range("F2").Select()
